$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), reusing the same header formatting (bold,
# bordered, centered) already applied to the other header cells like G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Corresponding value for row 2
$ws.Range("H2").Value = 1
